$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.306.95"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "2.371.10"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "506.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "130.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").Value = "2.380.57"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0986"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("E11").Value = "  -0.15%  "
$ws.Range("E12").Value = "  +6.57%  "
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("D14").Value = "2.792.71"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").Value = "56.283.96"
$ws.Range("E15").Value = "  -1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.79%  "
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "2.365.63"
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.98%  "
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "309.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.370"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("E27").Value = "  -2.55%  "
$ws.Range("E28").Value = "  -3.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("D30").Value = "0.0₃0711"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("E31").Value = "  -1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("E35").Value = "  -4.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.69"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.83%  "
$ws.Range("E39").Value = "  +1.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("E41").Value = "  -3.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "125.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.40%  "
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "238.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.87%  "
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.953"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.03%  "
